$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6 - this shifts the former rows 6..60 down to 7..61,
# matching the dimension change from A1:T60 to A1:T61.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new data record.
$ws.Range("A6").Value = 9
$ws.Range("B6").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C6").Value = "Metropolitana"
$ws.Range("D6").Value = 45063
$ws.Range("D6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E6").Value = 13
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100107
$ws.Range("H6").Value = "Otros"
$ws.Range("I6").Value = 100107001
$ws.Range("J6").Value = "Caqui"
$ws.Range("K6").Value = "Mankaki"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 420
$ws.Range("N6").Value = 13000
$ws.Range("O6").Value = 14000
$ws.Range("P6").Value = 13476
$ws.Range("Q6").Value = "`$/caja 16 kilos granel"
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 842
$ws.Range("T6").Value = 16
